$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Mapping Global Cyberterror Networks: An Empirical Study of Al-Qaeda and ISIS Cyberterrorism Events"
$ws.Range("C2").Value = "Claire Seungeun Lee, Kyung-Shick Choi, Ryan Shandler, Chris Kayser"
$ws.Range("D2").Value = "'2021"
$ws.Range("E2").Value = "10.1177/10439862211001606"
$ws.Range("G2").Value = "Journal of Contemporary Criminal Justice"

# Row 3
$ws.Range("B3").Value = "Indonesia’s Handling of Terrorists’ Cyber Activities: How Repressive Measures Still Fall Short"
$ws.Range("C3").Value = "Ali Abdullah Wibisono, Rachel Kumendong, Iwa Maulana"
$ws.Range("D3").Value = "'2025"
$ws.Range("E3").Value = "10.1177/23477970241298764"
$ws.Range("G3").Value = "Journal of Asian Security and International Affairs"

# Row 4
$ws.Range("B4").Value = "Shelf Life Prediction of Fresh Italian Pork Sausage Modified Atmosphere Packed"
$ws.Range("C4").Value = "E. Torrieri, F. Russo, R. Di Monaco, S. Cavella, F. Villani, F. Masi"
$ws.Range("D4").Value = "'2011"
$ws.Range("E4").Value = "10.1177/1082013210382328"
$ws.Range("G4").Value = "Food Science and Technology International"

# Row 5
$ws.Range("B5").Value = "A novel ensemble learning approach for fault detection of sensor data in cyber-physical system"
$ws.Range("C5").Value = "Ramesh Sneka Nandhini, Ramanathan Lakshmanan"
$ws.Range("D5").Value = "'2023"
$ws.Range("E5").Value = "10.3233/JIFS-235809"
$ws.Range("G5").Value = "Journal of Intelligent & Fuzzy Systems: Applications in Engineering and Technology"

# Row 6
$ws.Range("B6").Value = "Global versus Local Optimization in Redundancy Resolution of Robotic Manipulators"
$ws.Range("C6").Value = "Kazem Kazerounian, Zhaoyu Wang"
$ws.Range("D6").Value = "'1988"
$ws.Range("E6").Value = "10.1177/027836498800700501"
$ws.Range("G6").Value = "The International Journal of Robotics Research"

# Row 7
$ws.Range("B7").Value = "The code not taken: China, the United States, and the future of cyber espionage"
$ws.Range("C7").Value = "Adam Segal"
$ws.Range("D7").Value = "'2013"
$ws.Range("E7").Value = "10.1177/0096340213501344"
$ws.Range("G7").Value = "Bulletin of the Atomic Scientists"

# Row 8
$ws.Range("B8").Value = "Cyclones in cyberspace: Information shaping and denial in the 2008 Russia–Georgia war"
$ws.Range("C8").Value = "Ronald J. Deibert, Rafal Rohozinski, Masashi Crete-Nishihata"
$ws.Range("D8").Value = "'2012"
$ws.Range("E8").Value = "10.1177/0967010611431079"
$ws.Range("G8").Value = "Security Dialogue"
$ws.Range("L8").Value = 1

# Row 9
$ws.Range("B9").Value = "Offensive Cyber Operations and State Power: Lessons from Russia in Ukraine"
$ws.Range("C9").Value = "Alex S. Wilner, Gabriel Williams, Mattias Thuns-Rondeau, Nathanaël Beaulieu, Veronique Cossette-Sharkey"
$ws.Range("D9").Value = "'2024"
$ws.Range("E9").Value = "10.1177/00207020241234228"
$ws.Range("F9").Value = "Open Access"
$ws.Range("G9").Value = "International Journal: Canada’s Journal of Global Policy Analysis"

# Row 10
$ws.Range("B10").Value = "Framing cyber warfare: an analyst’s perspective"
$ws.Range("C10").Value = "Anthony Ween, Peter Dortmans, Nitin Thakur, Cayt Rowe"
$ws.Range("D10").Value = "'2019"
$ws.Range("E10").Value = "10.1177/1548512917725620"
$ws.Range("G10").Value = "The Journal of Defense Modeling and Simulation: Applications, Methodology, Technology"

# Row 11
$ws.Range("B11").Value = "Fighting in Cyberspace: Internet Access and the Substitutability of Cyber and Military Operations"
$ws.Range("C11").Value = "Nadiya Kostyuk, Erik Gartzke"
$ws.Range("D11").Value = "'2024"
$ws.Range("E11").Value = "10.1177/00220027231160993"
$ws.Range("F11").Value = "Restricted"
$ws.Range("G11").Value = "Journal of Conflict Resolution"
